# Cap nhat tien do
# Update progress-tracking dates for rows 6-8 (columns G/H) and
# refresh the sheet's active view (top-left scroll + selection).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 6: actual start/end dates
$ws.Range("G6").Value = 43383
$ws.Range("H6").Value = 43383

# Row 7: actual start/end dates
$ws.Range("G7").Value = 43384
$ws.Range("H7").Value = 43384

# Row 8: actual start/end dates
$ws.Range("G8").Value = 43386
$ws.Range("H8").Value = 43386

# Move the viewport so row 4 is at the top and select H9, matching
# the saved view state in the workbook.
$ws.Range("H9").Select()
$excel.ActiveWindow.ScrollRow = 4
$excel.ActiveWindow.ScrollColumn = 1
